$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "Source" block below the existing table (rows 6-8) ---
# Row 6: bold "Source:" label
$ws.Range("A6").Value = "Source:"
$ws.Range("A6").Font.Bold = $true

# Row 8 (URL) is written before row 7 (description) so that the new
# shared-string table entries land in the same order as the target file:
#   index 3 -> "Source:", index 4 -> URL, index 5 -> description
$ws.Range("A8").Value = "http://www.cdcr.ca.gov/Reports_Research/Offender_Information_Services_Branch/Quarterly/Strike1Archive.html"
$ws.Range("A7").Value = "California Department of Corrections and Rehabilitation: Second and Third Strike Inmate Population Report Archive"

# --- Misc view/print tweaks captured in the diff ---
# Selection moved to F22 before the file was saved.
[void]$ws.Range("F22").Select()

# Page orientation explicitly set to portrait (adds <pageSetup .../> to the sheet).
$ws.PageSetup.Orientation = 1
